# Applies the "2.1 first draft" update described in the commit message:
# - About sheet: source year corrected, source URL turned into a live
#   hyperlink, sheet no longer the tab shown on open
# - Data sheet: two corrected figures (typos in the WRI input data), becomes
#   the active/selected sheet, scrolled down to the bottom table
# - FoFObE sheet: entity categories renamed/expanded to match the new
#   (India 2.0) ownership breakdown, with a reworded, wrapped header

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("Data")
$wsFoF   = $wb.Worksheets.Item("FoFObE")

# ---------------------------------------------------------------------------
# About sheet
# ---------------------------------------------------------------------------
$wsAbout.Range("B5").Value = 2010

$wsAbout.Hyperlinks.Add(
    $wsAbout.Range("B7"),
    "http://www.fao.org/docrep/007/ae354e/AE354E10.htm#P587_26537",
    "P587_26537",
    "",
    "http://www.fao.org/docrep/007/ae354e/AE354E10.htm#P587_26537"
) | Out-Null

# ---------------------------------------------------------------------------
# Data sheet
# ---------------------------------------------------------------------------
$wsData.Range("B9").Value = 3170.51
$wsData.Range("B10").Value = 785.92

# ---------------------------------------------------------------------------
# FoFObE sheet
# ---------------------------------------------------------------------------
$wsFoF.Range("B1").Value = "Fraction of Forest Owned (dimensionless)"
$wsFoF.Rows.Item(1).RowHeight = 30
$wsFoF.Range("B1").HorizontalAlignment = -4152  # xlRight
$wsFoF.Range("B1").WrapText = $true

$wsFoF.Range("A3").Value = "nonenergy industries"

$wsFoF.Range("A4").Value = "labor and consumers"
$wsFoF.Range("B4").Value = 0
$wsFoF.Range("B4").NumberFormat = "0"

$wsFoF.Range("A5").Value = "foreign entities"
$wsFoF.Range("B5").Value = 0
$wsFoF.Range("B5").NumberFormat = "0"

$wsFoF.Range("A6").Value = "electricity suppliers"
$wsFoF.Range("B6").Value = 0
$wsFoF.Range("B6").NumberFormat = "0"

$wsFoF.Range("A7").Value = "coal suppliers"
$wsFoF.Range("B7").Value = 0
$wsFoF.Range("B7").NumberFormat = "0"

$wsFoF.Range("A8").Value = "natural gas and petroleum suppliers"
$wsFoF.Range("B8").Value = 0
$wsFoF.Range("B8").NumberFormat = "0"

$wsFoF.Range("A9").Value = "biomass and biofuel suppliers"
$wsFoF.Range("B9").Value = 0
$wsFoF.Range("B9").NumberFormat = "0"

$wsFoF.Range("A10").Value = "other energy suppliers"
$wsFoF.Range("B10").Value = 0
$wsFoF.Range("B10").NumberFormat = "0"

$wsFoF.Columns.Item(1).ColumnWidth = 32.67

# ---------------------------------------------------------------------------
# View state: per-sheet selections, and make Data the sheet shown on open
# ---------------------------------------------------------------------------
$wsAbout.Select()
$wsAbout.Range("B10").Select() | Out-Null

$wsFoF.Select()
$wsFoF.Range("E12").Select() | Out-Null

$wsData.Select()
$wsData.Range("E26").Select() | Out-Null

$wb.Save()
